$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.008.64"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.801.95"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.06%  "
$ws.Range("D7").Value = "3.801.94"
$ws.Range("E7").Value = "  +2.01%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.458"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000244"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").Value = "4.437.73"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "3.808.23"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("D17").Value = "69.161.60"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.52%  "
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "486.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.717"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000158"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.54%  "
$ws.Range("D33").Value = "3.960.83"
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "3.748.38"
$ws.Range("E35").Value = "  +2.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.106"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.139"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.317"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "436.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("D48").Value = "2.822.31"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0352"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.19%  "
